{"js": "// Rewrite the Italy/NATO paragraph:\n//   1. \"...NATO however had a complicated...\"   -> \"...NATO however have had a complicated...\"\n//   2. \"...different presidents sat, they...\"   -> \"...different presidents led the country, they...\"\n//   3. \"...far to the Left.  Internally however the tumultuous...\" -> \"...far to the Left.  Internally the tumultuous...\"\n// and drop the stray \"_GoBack\" bookmark left over from the previous edit session.\n\nconst body = context.document.body;\n\nasync function findAndReplace(findText, replaceText) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nawait findAndReplace(\"however had\", \"however have had\");\nawait findAndReplace(\"presidents sat,\", \"presidents led the country,\");\nawait findAndReplace(\"Internally however the\", \"Internally the\");\n\n// Remove the leftover \"_GoBack\" bookmark (Word inserts this automatically at the\n// last edit point; it is not meant to persist once the edit is finalized).\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Rewrite the Italy/NATO paragraph:\n#   1. \"...NATO however had a complicated...\"   -> \"...NATO however have had a complicated...\"\n#   2. \"...different presidents sat, they...\"   -> \"...different presidents led the country, they...\"\n#   3. \"...far to the Left.  Internally however the tumultuous...\" -> \"...far to the Left.  Internally the tumultuous...\"\n# and drop the stray \"_GoBack\" bookmark left over from the previous edit session.\n\n$d = $word.ActiveDocument\n\n$rng1 = $d.Content\n$rng1.Find.Execute(\"however had\", $false, $false, $false, $false, $false, $true, 1, $false, \"however have had\", 2) | Out-Null\n\n$rng2 = $d.Content\n$rng2.Find.Execute(\"presidents sat,\", $false, $false, $false, $false, $false, $true, 1, $false, \"presidents led the country,\", 2) | Out-Null\n\n$rng3 = $d.Content\n$rng3.Find.Execute(\"Internally however the\", $false, $false, $false, $false, $false, $true, 1, $false, \"Internally the\", 2) | Out-Null\n\n# Remove the leftover \"_GoBack\" bookmark (Word inserts this automatically at the\n# last edit point; it is not meant to persist once the edit is finalized).\n$d.Bookmarks(\"_GoBack\").Delete()\n"}
